# Versão 4 = Consegue fazer o download da guia itens e anexo ( concomitantemente )
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: requisição 8167172 -> 8167085, status Aprovado -> Cancelado
$ws.Range("A2").Value = 8167085
$ws.Range("C2").Value = "Cancelado"

# Row 3: requisição 8167085 -> 8166505, status Cancelado -> Aprovado
$ws.Range("A3").Value = 8166505
$ws.Range("C3").Value = "Aprovado"

# Preserve the number-formatting (Consolas, vertically centered) that lived
# on D5 by copying it onto D4 before the old row is cleared out.
$ws.Range("D5").Copy()
$ws.Range("D4").PasteSpecial(-4122)

# New row 4: requisição 8167172, status Aprovado, D4 stays an empty
# formatted cell (format only, no value)
$ws.Range("A4").Value = 8167172
$ws.Range("C4").Value = "Aprovado"
$ws.Range("D4").ClearContents()

# Remove the old trailing row (previously row 5, columns D:F) that is no
# longer part of the table
$ws.Range("D5").Clear()
$ws.Range("E5").Clear()
$ws.Range("F5").Clear()

$ws.Range("C3").Select()
